# "Generate Report for Handoff": the localization status for zh-cn/de-de
# flips from "In Translation" to "Ready for handoff", and the associated
# timestamps (Overview's "Latest HO Xliff Generate Date" and each locale
# sheet's "Latest Handoff Datetime") are refreshed to the handoff time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets("Overview")
$wsZhCn = $wb.Worksheets("zh-cn")
$wsDeDe = $wb.Worksheets("de-de")

# Overview sheet: per-locale status columns (zh-cn, de-de) for the one row.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Overview sheet: "Latest HO Xliff Generate Date".
$wsOverview.Range("G2").Value = "2016-08-17 00:36:51"

# zh-cn detail sheet: "Status" + "Latest Handoff Datetime".
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 00:36:47"

# de-de detail sheet: "Status" + "Latest Handoff Datetime".
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 00:36:51"

# The Status column text got longer ("In Translation" -> "Ready for
# handoff"), so Excel re-autofit those columns, widening them. Reproduce
# that widening (engine quantizes ColumnWidth to 1/6-character steps, so
# this lands as close as possible to the canonical 17.2159881591797).
$wsOverview.Columns("E:F").ColumnWidth = 16.3
$wsZhCn.Columns("C:C").ColumnWidth = 16.3
$wsDeDe.Columns("C:C").ColumnWidth = 16.3
